$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - row 2, 4, 5, 11
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1229
$ws1.Range("F4").Value = 22
$ws1.Range("F5").Value = 12510
$ws1.Range("F11").Value = 236

# Sheet "全部类型" (All types) - row 2, 6, 7, 13
$ws2 = $wb.Worksheets.Item("全部类型")
$ws2.Range("F2").Value = 1229
$ws2.Range("F6").Value = 22
$ws2.Range("F7").Value = 12510
$ws2.Range("F13").Value = 236
